$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$games = @(
    @(16, "LA_SEA", 44.5, -1.5),
    @(16, "PHI_WAS", 48.5, -1.5),
    @(16, "GB_CHI", 45.5, -1.5),
    @(16, "BUF_CLE", 43.5, -8.5),
    @(16, "NYJ_NO", 41.5, 1.5),
    @(16, "MIN_NYG", 42.5, -2.5),
    @(16, "TB_CAR", 47.5, -3.5),
    @(16, "KC_TEN", 44.5, -7),
    @(16, "LAC_DAL", 46.5, -1.5),
    @(16, "CIN_MIA", 48.5, -1.5),
    @(16, "ATL_ARI", 46.5, 3.5),
    @(16, "JAX_DEN", 45.5, 5.5),
    @(16, "PIT_DET", 45.5, 6),
    @(16, "LV_HOU", 43.5, 4.5),
    @(16, "NE_BAL", 46.5, 8.5),
    @(16, "SF_IND", 46.5, -3)
)

$startRow = 210
for ($i = 0; $i -lt $games.Count; $i++) {
    $row = $startRow + $i
    $g = $games[$i]
    $ws.Cells.Item($row, 1).Value = $g[0]
    $ws.Cells.Item($row, 2).Value = $g[1]
    $ws.Cells.Item($row, 3).Value = $g[2]
    $ws.Cells.Item($row, 4).Value = $g[3]
}

$null = $ws.Range("D210").Select()
